$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row: "..._old" -> "..._FV2410" and "..._new" -> "..._FV2504"
$oldHeaders = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old"
)
$newHeaders = @(
    "Segmentname_new",
    "Segmentgruppe_new",
    "Segment_new",
    "Datenelement_new",
    "Segment ID_new",
    "Code_new",
    "Qualifier_new",
    "Beschreibung_new",
    "Bedingungsausdruck_new",
    "Bedingung_new"
)

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $oldHeaders[$i].Replace("_old", "_FV2410")
}

# Column K (11) holds "diff" and stays untouched.

for ($i = 0; $i -lt $newHeaders.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $newHeaders[$i].Replace("_new", "_FV2504")
}

# 2. Turn the used range A1:U90 into an Excel Table (ListObject) with an autofilter.
$tableRange = $ws.Range("A1:U90")
$lo = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"

# 3. Freeze the header row (split/freeze after row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
